$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.515.83"
$ws.Range("E2").Value = "  +5.21%  "
$ws.Range("D3").Value = "1.724.75"
$ws.Range("E3").Value = "  +4.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "225.87"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").Value = "0.5373"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").Value = "21.79"
$ws.Range("E10").Value = "  +5.93%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "4.638"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.711.01"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "1.961.24"
$ws.Range("E14").Value = "  +4.10%  "
$ws.Range("D15").Value = "0.5875"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "0.0₅8293"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "68.12"
$ws.Range("E17").Value = "  +3.93%  "
$ws.Range("D18").Value = "27.531.54"
$ws.Range("E18").Value = "  +5.30%  "
$ws.Range("D19").Value = "223.56"
$ws.Range("E19").Value = "  +15.77%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "4.745"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "6.109"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "148.25"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "1.694"
$ws.Range("E26").Value = "  +11.56%  "
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("D28").Value = "7.417"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "16.71"
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").Value = "0.05553"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "1.303"
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("D32").Value = "3.557"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "3.469"
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("D34").Value = "1.664"
$ws.Range("D35").Value = "0.9610"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "2.445"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").Value = "2.817"
$ws.Range("D38").Value = "0.5956"
$ws.Range("E38").Value = "  +4.48%  "
$ws.Range("D39").Value = "0.01651"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("D40").Value = "5.874"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "1.060.51"
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("D42").Value = "0.8571"
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "101.66"
$ws.Range("D45").Value = "1.867.07"
$ws.Range("E45").Value = "  +4.00%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").Value = "59.09"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").Value = "8.202"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "0.4440"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "0.05272"
$ws.Range("E51").Value = "  +1.36%  "
